# fix(publipostage): Add space before ":"
#
# - add a space before the ":" in the "statut_name" column (B) values
# - add two new EudraCT numbers (column D) for rows 9 and 11
# - rotate the NCTId/title/acronym/intervention_type (C/G/H/I) content of
#   rows 11-13 by one position
# - correct a few intervention_type (column I) values
# - add a new intervention_type value for row 16 (column I)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: "statut_name" -- add a space before the colon ---------------
$statut4 = "4 : pas de résultats postés ni publiés"
$statut1 = "1 : résultats postés ou publiés dans les 12 mois"

$ws.Range("B2").Value  = $statut4
$ws.Range("B3").Value  = $statut4
$ws.Range("B4").Value  = $statut4
$ws.Range("B5").Value  = $statut4
$ws.Range("B6").Value  = $statut4
$ws.Range("B7").Value  = $statut4
$ws.Range("B8").Value  = $statut4
$ws.Range("B9").Value  = $statut1
$ws.Range("B10").Value = $statut4
$ws.Range("B11").Value = $statut4
$ws.Range("B12").Value = $statut4
$ws.Range("B13").Value = $statut4
$ws.Range("B14").Value = $statut4
$ws.Range("B15").Value = $statut4
$ws.Range("B16").Value = $statut4

# --- New EudraCT numbers (column D) -----------------------------------------
$ws.Range("D9").Value  = "2020-001570-30"
$ws.Range("D11").Value = "2013-002056-33"

# --- Rows 11-13: rotate NCTId / title / acronym / intervention_type --------
$ws.Range("C11").Value = "NCT02235012"
$ws.Range("G11").Value = "Cognitive Biases in Decision Making in a Pharmacological Model of Psychosis : a Study in Healthy Humans Recieving Low Dose Anesthetic, Ketamine Versus Placebo"
$ws.Range("H11").Value = "KETABI"
$ws.Range("I11").Value = "DRUG"

$ws.Range("C12").Value = "NCT02841098"
$ws.Range("G12").Value = """ Endarterectomy Combined With Optimal Medical Therapy Versus Optimal Medical Therapy Alone in Patients With Asymptomatic Severe Atherosclerotic Carotid Artery Stenosis at Higher-than-average Risk of Ipsilateral Stroke """
$ws.Range("H12").Value = "ACTRIS"
$ws.Range("I12").Value = "OTHER"

$ws.Range("C13").Value = "NCT02476435"
$ws.Range("G13").Value = "Depersonalization Disorder: Therapeutic Effect of Neuronavigated Repetitive Transcranial Stimulation of Right Angular Gyrus."
$ws.Range("H13").Value = "PERSONA"
$ws.Range("I13").Value = "DEVICE"

# --- Other intervention_type (column I) corrections -------------------------
$ws.Range("I2").Value  = "OTHER"
$ws.Range("I14").Value = "OTHER"
$ws.Range("I16").Value = "DRUG (presumed)"
